$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 to the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-10
$values = @(
    @(7, 8),
    @(7, 8),
    @(7, 8),
    @(6, 8),
    @(7, 8),
    @(7, 9),
    @(7, 9),
    @(4, 6),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
